$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the three "CA ..." header labels to "CA Radius ..." (LMT-D misalignment
# fix: align header wording with the LMT-A sheet / same channels).
$ws.Range("G1").Value = "CA Radius (Thorlabs 90% visible portion of lens)"
$ws.Range("H1").Value = "CA Radius req"
$ws.Range("J1").Value = "CA Radius (zmx model from Thorlabs)"

# Match the author's final selection state (merged cell J2:J4).
$ws.Range("J2:J4").Select()
